$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 665 (old rows 665-765 shift down to 667-767).
$ws.Rows.Item(665).Insert()
$ws.Rows.Item(665).Insert()

# New row 665
$ws.Range("A665").Value = 3
$ws.Range("B665").Value = "Femacal de La Calera"
$ws.Range("C665").Value = "Coquimbo"
$ws.Range("D665").Value = 45077
$ws.Range("E665").Value = 5
$ws.Range("F665").Value = 100112003
$ws.Range("G665").Value = "Ajo"
$ws.Range("H665").Value = "Chino"
$ws.Range("I665").Value = "Primera"
$ws.Range("J665").Value = 80
$ws.Range("K665").Value = 15500
$ws.Range("L665").Value = 16000
$ws.Range("M665").Value = 15750
$ws.Range("N665").Value = "$/caja 10 kilos"
$ws.Range("O665").Value = "China"
$ws.Range("P665").Value = 1575
$ws.Range("Q665").Value = 10
$ws.Range("R665").Value = "Hortaliza"

# New row 666
$ws.Range("A666").Value = 3
$ws.Range("B666").Value = "Femacal de La Calera"
$ws.Range("C666").Value = "Coquimbo"
$ws.Range("D666").Value = 45077
$ws.Range("E666").Value = 5
$ws.Range("F666").Value = 100112003
$ws.Range("G666").Value = "Ajo"
$ws.Range("H666").Value = "Chino"
$ws.Range("I666").Value = "Primera"
$ws.Range("J666").Value = 45
$ws.Range("K666").Value = 17000
$ws.Range("L666").Value = 17000
$ws.Range("M666").Value = 17000
$ws.Range("N666").Value = "$/malla 10 kilos"
$ws.Range("O666").Value = "China"
$ws.Range("P666").Value = 1700
$ws.Range("Q666").Value = 10
$ws.Range("R666").Value = "Hortaliza"
